# Inserts a new weekly price-report row above the current row 25, shifting
# all subsequent rows (old 25-55) down to (26-56), and populates the new
# row 25 with this week's data for "Poroto verde" / "Sin especificar" / "Primera".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 25-55 down by inserting a new row at 25.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C25").Value = "Arica y Parinacota"
$ws.Range("D25").Value = 44651
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = 100112031
$ws.Range("G25").Value = "Poroto verde"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 1500
$ws.Range("K25").Value = 600
$ws.Range("L25").Value = 700
$ws.Range("M25").Value = 650
$ws.Range("N25").Value = "$/kilo"
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("P25").Value = 650
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = "Hortaliza"
